# Auto-generated Excel COM-interop script
# Applies 2025-11-03 YTD violent crime data updates (column L = 2025) across
# the Citywide Totals, By Neighborhood summary, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5658
$ws.Range("L3").Value = 6143
$ws.Range("L4").Value = 1511
$ws.Range("L5").Value = 364
$ws.Range("L6").Value = 5051
$ws.Range("L7").Value = 18727

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L4").Value = 68
$ws.Range("L8").Value = 1243
$ws.Range("L9").Value = 109
$ws.Range("L11").Value = 307
$ws.Range("L18").Value = 127
$ws.Range("L19").Value = 513
$ws.Range("L20").Value = 466
$ws.Range("L24").Value = 49
$ws.Range("L27").Value = 165
$ws.Range("L29").Value = 1056
$ws.Range("L31").Value = 184
$ws.Range("L33").Value = 852
$ws.Range("L36").Value = 238
$ws.Range("L37").Value = 715
$ws.Range("L42").Value = 606
$ws.Range("L44").Value = 128
$ws.Range("L48").Value = 244
$ws.Range("L51").Value = 233
$ws.Range("L52").Value = 388
$ws.Range("L54").Value = 414
$ws.Range("L55").Value = 196
$ws.Range("L57").Value = 66
$ws.Range("L63").Value = 57
$ws.Range("L67").Value = 650
$ws.Range("L69").Value = 42
$ws.Range("L72").Value = 79
$ws.Range("L76").Value = 289
$ws.Range("L77").Value = 125
$ws.Range("L79").Value = 509
$ws.Range("L80").Value = 61
$ws.Range("L83").Value = 413
$ws.Range("L84").Value = 181
$ws.Range("L85").Value = 929
$ws.Range("L90").Value = 193
$ws.Range("L95").Value = 264
$ws.Range("L96").Value = 213
$ws.Range("L98").Value = 100
$ws.Range("L100").Value = 32
$ws.Range("L101").Value = 18727

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L6").Value = 64
$ws.Range("L7").Value = 213

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L6").Value = 73
$ws.Range("L7").Value = 307

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 284
$ws.Range("L3").Value = 381
$ws.Range("L7").Value = 929

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 124
$ws.Range("L7").Value = 388

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 369
$ws.Range("L3").Value = 438
$ws.Range("L4").Value = 88
$ws.Range("L7").Value = 1243

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 131
$ws.Range("L3").Value = 168
$ws.Range("L7").Value = 413

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 241
$ws.Range("L7").Value = 852

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 99
$ws.Range("L7").Value = 264

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 214
$ws.Range("L6").Value = 193
$ws.Range("L7").Value = 715

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 184

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 251
$ws.Range("L4").Value = 43
$ws.Range("L6").Value = 150
$ws.Range("L7").Value = 650

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 59
$ws.Range("L7").Value = 181

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 200
$ws.Range("L7").Value = 414

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 311
$ws.Range("L3").Value = 407
$ws.Range("L7").Value = 1056

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 63
$ws.Range("L7").Value = 244

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 182
$ws.Range("L3").Value = 159
$ws.Range("L7").Value = 513

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 50
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 128

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 131
$ws.Range("L7").Value = 289

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 165
$ws.Range("L7").Value = 606

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 67
$ws.Range("L7").Value = 196

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 166
$ws.Range("L3").Value = 165
$ws.Range("L6").Value = 130
$ws.Range("L7").Value = 509

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 146
$ws.Range("L3").Value = 158
$ws.Range("L7").Value = 466

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 47
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 238

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L2").Value = 33
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 165

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 193

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 73
$ws.Range("L7").Value = 233

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L3").Value = 41
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L4").Value = 7
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L5").Value = 2
$ws.Range("L7").Value = 68
